$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Settings")
$ws2 = $wb.Worksheets.Item("Constants")

# Add new Orchestrator REST API settings rows (4-8) on the Settings sheet.
$ws1.Range("A4").Value = "Orch_tenancyName"
$ws1.Range("B4").Value = "fantastic"

$ws1.Range("A5").Value = "Orch_env"
$ws1.Range("B5").Value = "mihhdusENV"

$ws1.Range("A6").Value = "Orch_orchestratorURL"
$ws1.Range("B6").Value = "https://demo.uipath.com"

$ws1.Range("A7").Value = "Orch_userNameOrEmailAddress"
$ws1.Range("B7").Value = "admin"

$ws1.Range("A8").Value = "Orch_password"
$ws1.Range("B8").Value = "123qwe"

# Update the queue name used on row 2 (written last so it lands at the
# end of the shared-string table, matching the authored workbook).
$ws1.Range("B2").Value = "KibanaDemoQueue"

# Update the selection on the Constants sheet (now inactive) before
# switching the active tab back to Settings.
[void]$ws2.Range("A24").Select()

[void]$ws1.Activate()
[void]$ws1.Range("B2").Select()
